function Set-CellText {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    $cell = $Sheet.Range($Address)
    if ($Text.Trim() -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        # Value would otherwise be auto-parsed as a number by Excel;
        # force text formatting, assign, then strip the format marker
        # so the cell keeps its original (unstyled) appearance.
        $cell.NumberFormat = "@"
        $cell.Value = $Text
        $cell.ClearFormats()
    } else {
        $cell.Value = $Text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" '29.025.44'
Set-CellText $ws "E2" '  -0.53%  '
Set-CellText $ws "D3" '1.829.61'
Set-CellText $ws "E3" '  -0.32%  '
Set-CellText $ws "D4" '0.9987'
Set-CellText $ws "E4" '  -0.06%  '
Set-CellText $ws "D5" '241.18'
Set-CellText $ws "E5" '  -0.31%  '
Set-CellText $ws "D6" '0.6250'
Set-CellText $ws "E6" '  -5.02%  '
Set-CellText $ws "D7" '0.9999'
Set-CellText $ws "E7" '  +0.00%  '
Set-CellText $ws "D8" '0.07570'
Set-CellText $ws "E8" '  +1.95%  '
Set-CellText $ws "D9" '44.90'
Set-CellText $ws "E9" '  +7.33%  '
Set-CellText $ws "D10" '0.2910'
Set-CellText $ws "E10" '  -0.59%  '
Set-CellText $ws "D11" '22.77'
Set-CellText $ws "E11" '  -0.77%  '
Set-CellText $ws "D12" '0.07638'
Set-CellText $ws "E12" '  -1.70%  '
Set-CellText $ws "D13" '1.831.85'
Set-CellText $ws "E13" '  -2.35%  '
Set-CellText $ws "D14" '4.952'
Set-CellText $ws "E14" '  -0.61%  '
Set-CellText $ws "D15" '0.6646'
Set-CellText $ws "E15" '  -0.12%  '
Set-CellText $ws "D16" '82.23'
Set-CellText $ws "E16" '  -0.70%  '
Set-CellText $ws "D17" '0.000009111'
Set-CellText $ws "E17" '  +5.97%  '
Set-CellText $ws "E18" '  -2.11%  '
Set-CellText $ws "D19" '28.924.24'
Set-CellText $ws "E19" '  -0.95%  '
Set-CellText $ws "D20" '224.68'
Set-CellText $ws "E20" '  -0.87%  '
Set-CellText $ws "E21" '  -1.10%  '
Set-CellText $ws "D22" '1.000'
Set-CellText $ws "E22" '  +0.03%  '
Set-CellText $ws "D23" '7.191'
Set-CellText $ws "E23" '  +0.91%  '
Set-CellText $ws "D24" '1.000'
Set-CellText $ws "E24" '  +0.01%  '
Set-CellText $ws "D25" '159.78'
Set-CellText $ws "E25" '  +0.39%  '
Set-CellText $ws "D26" '8.406'
Set-CellText $ws "E26" '  -2.28%  '
Set-CellText $ws "D27" '0.1360'
Set-CellText $ws "E27" '  -2.55%  '
Set-CellText $ws "D28" '17.80'
Set-CellText $ws "E28" '  -0.75%  '
Set-CellText $ws "D29" '1.497'
Set-CellText $ws "E29" '  -1.14%  '
Set-CellText $ws "E30" '  -0.49%  '
Set-CellText $ws "D31" '4.044'
Set-CellText $ws "E31" '  -1.69%  '
Set-CellText $ws "E32" '  +0.55%  '
Set-CellText $ws "D33" '0.05198'
Set-CellText $ws "E33" '  -1.32%  '
Set-CellText $ws "D34" '1.843'
Set-CellText $ws "E34" '  -1.23%  '
Set-CellText $ws "D35" '1.152'
Set-CellText $ws "E35" '  +0.51%  '
Set-CellText $ws "D36" '0.7305'
Set-CellText $ws "E36" '  -1.24%  '
Set-CellText $ws "E37" '  -1.60%  '
Set-CellText $ws "D38" '1.277.77'
Set-CellText $ws "E38" '  -2.19%  '
Set-CellText $ws "D39" '2.758'
Set-CellText $ws "E39" '  +0.93%  '
Set-CellText $ws "D40" '0.01789'
Set-CellText $ws "E40" '  -0.40%  '
Set-CellText $ws "D41" '6.403'
Set-CellText $ws "E41" '  +5.90%  '
Set-CellText $ws "D42" '0.8887'
Set-CellText $ws "E42" '  -3.48%  '
Set-CellText $ws "D44" '101.59'
Set-CellText $ws "E44" '  -0.84%  '
Set-CellText $ws "D45" '1.978.51'
Set-CellText $ws "E45" '  -2.42%  '
Set-CellText $ws "D46" '0.5109'
Set-CellText $ws "E46" '  -0.65%  '
Set-CellText $ws "D47" '63.69'
Set-CellText $ws "E47" '  +0.17%  '
Set-CellText $ws "D48" '0.00000000119'
Set-CellText $ws "E48" '  -0.98%  '
Set-CellText $ws "B49" 'TheSandbox'
Set-CellText $ws "C49" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws "D49" '0.3976'
Set-CellText $ws "E49" '  -0.64%  '
Set-CellText $ws "B50" 'XinFinNetwork'
Set-CellText $ws "C50" 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-CellText $ws "D50" '0.07269'
Set-CellText $ws "E50" '  -15.95%  '
Set-CellText $ws "D51" '8.860'
Set-CellText $ws "E51" '  +1.66%  '
